# lesson 4 exercises completed
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B39 has no prior formatting (unlike the rest of column B in this block),
# so copy the "fill-in" cell style from a neighboring cell before setting it.
$ws.Range("B40").Copy()
$ws.Range("B39").PasteSpecial(-4122)

# Fill in the Ukrainian sentences (column B) and their English translations
# (column D) for LESSON 4 (rows 39-48), in the same order the values were
# entered so the shared-string table ends up in the same sequence.
$ws.Range("B39").Value2 = "Ми зараз почнемо збори"
$ws.Range("B40").Value2 = "Вона подвонить тобі на перерві"
$ws.Range("B41").Value2 = "Ти багато вивчиш у цьому курсі"
$ws.Range("B42").Value2 = "Вони прибудуть біля 6 вечора"
$ws.Range("B43").Value2 = "Джон залишиться на роботі допізна"
$ws.Range("B44").Value2 = "Він сформує нову команду"
$ws.Range("B45").Value2 = "Я приготую м'со на вечерю"
$ws.Range("B47").Value2 = "Повітря буде відчуватися холоднішим"
$ws.Range("D48").Value2 = "Wt will spend the day together"
$ws.Range("B48").Value2 = "Ми проведемо день разом"
$ws.Range("D47").Value2 = "The air will feel cooler"
$ws.Range("D39").Value2 = "We will start the meeting now"
$ws.Range("D40").Value2 = "She will call you during the break"
$ws.Range("D41").Value2 = "You will learn a lot in this course"
$ws.Range("D42").Value2 = "They will arrive around 6p.m."
$ws.Range("D43").Value2 = "John will stay late at work"
$ws.Range("D44").Value2 = "He will form a new team"
$ws.Range("D45").Value2 = "I will cook meat for dinner"
$ws.Range("D46").Value2 = "She will visit this place again"
$ws.Range("B46").Value2 = "Вона відвідає це місце знову"

# Update the window/view state to reflect where the author was working.
$ws.Application.ActiveWindow.ScrollRow = 35
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("D42").Select()
